$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3 data: Aanand (name), email, SL=2, CL=6
$ws.Range("A3").Value = "Aanand "
$ws.Range("B3").Value = "anand@monetnetworks.com"
$ws.Range("C3").Value = 2
$ws.Range("E3").Value = 6

# Row 2 gains a CL value
$ws.Range("E2").Value = 1

# Update selection to E3
$ws.Range("E3").Select()
